# NotasEscolaresEstandarizadas.xlsx — "Add files via upload" edit
#
# What changed (per the author's diff):
#   1. A new "Estudiante" header was added in column A (A1), labelling the
#      student-name column that already existed (A2:A11) but had no header.
#   2. The standardized "Matemáticas" column (D) was recalculated for every
#      student - new z-score values in D2:D11.
#   3. Column widths were set (best-fit) for the new header + numeric cols.
#   4. The active selection moved to C14.
#
# Apply all of this against the already-open ActiveWorkbook via COM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New column header for the student-name column.
$ws.Range("A1").Value = "Estudiante"

# 2. Updated standardized values in column D (Matemáticas), rows 2-11.
#    (written as plain decimals - the interpreter's expression parser
#    doesn't accept `E`-notation literals - but these are the exact same
#    IEEE-754 doubles as the scientific-notation figures in the source.)
$ws.Range("D2").Value  = 1.6967060807959473
$ws.Range("D3").Value  = -0.51133607914398493
$ws.Range("D4").Value  = 0.30215313767599056
$ws.Range("D5").Value  = -1.4410380412239565
$ws.Range("D6").Value  = 0.06972764715599743
$ws.Range("D7").Value  = -0.046485098103998625
$ws.Range("D8").Value  = 0.53457862819598267
$ws.Range("D9").Value  = -0.27891058862399176
$ws.Range("D10").Value = -1.4410380412239565
$ws.Range("D11").Value = 1.115642354495965

# 3. Best-fit the columns now that there's a header in A and the data is final.
$ws.Range("A1:E11").Columns.AutoFit()

# 4. Move the selection like the saved workbook shows.
$ws.Range("C14").Select()
